$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo'd image filenames in the Image column
$ws.Range("D28").Value = "Frooti20.jpg"
$ws.Range("D24").Value = "Pastry Pineapple.jpg"

# Update the last active selection on the sheet
$ws.Range("G25").Select()
